$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing thread_id values in D2 and D3
$ws.Range("D2").Value = 'profile-thread-d8c18ce5-351c-45ae-88f2-b7aa77bbec78'
$ws.Range("D3").Value = 'profile-thread-ffe4e96f-25b0-4dec-bd45-5ada36708478'

# Append a new error row (row 4)
$ws.Range("A4").Value = 'https://www.uidaho.edu/cals/agricultural-economics-and-rural-sociology/our-people/andres-trujillo-barrera'
$ws.Range("B4").Value = 'LLM output failed validation/parsing after successful call'
$ws.Range("C4").Value = '{''exception_type'': ''ValidationError'', ''message'': "3 validation errors for ProfileData\ndegrees.0\n  Input should be a valid string [type=string_type, input_value={''degree'': ''Ph.D.'', ''inst...hampaign'', ''year'': 2013}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type\ndegrees.1\n  Input should be a valid string [type=string_type, input_value={''degree'': ''M.S.'', ''insti...iversity'', ''year'': 2006}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type\ndegrees.2\n  Input should be a valid string [type=string_type, input_value={''degree'': ''B.S.'', ''insti...a Sabana'', ''year'': 2001}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type", ''llm_content'': ''```json\n{\n  "first_name": "Andres",\n  "middle_name": null,\n  "last_name": "Trujillo-Barrera",\n  "title": "Associate Professor & Director, Agricultural Commodity Risk Management Program",\n  "office": "Ag Science, Room 37",\n  "phone": "208-885-1151",\n  "email": "aatrujillo@uidaho.edu",\n  "college_unit": "College of Agricultural and Life Sciences",\n  "department_division": "Department of Agricultural Economics and Rural Sociology",\n  "degrees": [\n    {\n      "degree": "Ph.D.",\n      "institution": "University of Illinois at Urbana-Champaign",\n      "year": 2013\n    },\n    {\n      "degree": "M.S.",\n      "institution": "Illinois State University",\n      "year": 2006\n    },\n    {\n      "degree": "B.S.",\n      "institution": "Universidad de la Sabana",\n      "year": 2001\n    }\n  ],\n  "research_focus_areas": [\n    "A better understanding of volatility and its effect on agribusiness supply chains (marketing strategies, performance, industrial organization)",\n    "Evaluation of risk management instruments (futures, options, insurance, market microstructure)",\n    "Application of new tools to measure and evaluate increased uncertainty by using a combination of marketing, economics and finance theory with quantitative research methods"\n  ],\n  "photo_url": null\n}\n```'', ''traceback'': ''Traceback (most recent call last):\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/src/nodes.py", line 474, in extract_data\n    extracted_profile = ProfileData(**target_obj)\n  File "/Users/colesummers/Documents/GitHub/profile-extractor/.venv/lib/python3.9/site-packages/pydantic/main.py", line 243, in __init__\n    validated_self = self.__pydantic_validator__.validate_python(data, self_instance=self)\npydantic_core._pydantic_core.ValidationError: 3 validation errors for ProfileData\ndegrees.0\n  Input should be a valid string [type=string_type, input_value={''degree'': ''Ph.D.'', ''inst...hampaign'', ''year'': 2013}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type\ndegrees.1\n  Input should be a valid string [type=string_type, input_value={''degree'': ''M.S.'', ''insti...iversity'', ''year'': 2006}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type\ndegrees.2\n  Input should be a valid string [type=string_type, input_value={''degree'': ''B.S.'', ''insti...a Sabana'', ''year'': 2001}, input_type=dict]\n    For further information visit https://errors.pydantic.dev/2.11/v/string_type\n''}'
$ws.Range("D4").Value = 'profile-thread-047df88e-28c8-4586-b00b-8ca4b7ea5322'

Write-Host "Row 4 added; D2/D3 updated"
